$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2232.375
$ws.Range("J62").Value = 3339.75
$ws.Range("L62").Value = 3339.75
$ws.Range("N62").Value = -4587.75
$ws.Range("H65").Value = 2232.375
$ws.Range("J65").Value = 3339.75
$ws.Range("L65").Value = 16698.75
$ws.Range("N65").Value = -22938.75
$ws.Range("H68").Value = 29859
$ws.Range("J68").Value = 29859
$ws.Range("L68").Value = 29859
$ws.Range("N68").Value = -31357
$ws.Range("H71").Value = 29859
$ws.Range("J71").Value = 29859
$ws.Range("L71").Value = 89577
$ws.Range("N71").Value = -97065
$ws.Range("H75").Value = 27078.75
$ws.Range("J75").Value = 27078.75
$ws.Range("L75").Value = 27078.75
$ws.Range("N75").Value = -28950.75
$ws.Range("H78").Value = 27078.75
$ws.Range("J78").Value = 27078.75
$ws.Range("L78").Value = 81236.25
$ws.Range("N78").Value = -90596.25
$ws.Range("H98").Value = 3612.6
$ws.Range("I98").Value = 3820.6428
$ws.Range("J98").Value = 700
$ws.Range("K98").Value = 3820.6428
$ws.Range("L98").Value = 700
$ws.Range("M98").Value = -2322.6428
$ws.Range("N98").Value = -3696
$ws.Range("H113").Value = 3926.72
$ws.Range("I113").Value = 3377.7646
$ws.Range("J113").Value = 5093.25
$ws.Range("K113").Value = 3377.7646
$ws.Range("L113").Value = 5093.25
$ws.Range("M113").Value = -123.7646
$ws.Range("N113").Value = -11601.25
$ws.Range("H116").Value = 3681.318
$ws.Range("I116").Value = 3587
$ws.Range("J116").Value = 3759.9167
$ws.Range("K116").Value = 3587
$ws.Range("L116").Value = 3759.9167
$ws.Range("M116").Value = -145
$ws.Range("N116").Value = -10643.9167
$ws.Range("H122").Value = 3612.6
$ws.Range("I122").Value = 3820.6428
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 11461.9284
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -9011.9284
$ws.Range("N122").Value = -7000

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3040.7727
$ws.Range("I61").Value = 2798.3333
$ws.Range("K61").Value = 2798.3333
$ws.Range("M61").Value = -2586.3333
$ws.Range("H74").Value = 1181.4166
$ws.Range("I74").Value = 826
$ws.Range("K74").Value = 826
$ws.Range("M74").Value = 48
$ws.Range("H77").Value = 1181.4166
$ws.Range("I77").Value = 826
$ws.Range("K77").Value = 4130
$ws.Range("M77").Value = 238
$ws.Range("H86").Value = 25333.334
$ws.Range("H89").Value = 25333.334
$ws.Range("H132").Value = 2172.6482
$ws.Range("I132").Value = 1540.814
$ws.Range("K132").Value = 4622.442
$ws.Range("M132").Value = -2092.442
$ws.Range("H136").Value = 3040.7727
$ws.Range("I136").Value = 2798.3333
$ws.Range("K136").Value = 8394.999899999999
$ws.Range("M136").Value = -5844.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 34992.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 34992.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 34992.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -36364.5
$ws.Range("H65").Value = 34992.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 34992.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 104977.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -111841.5
$ws.Range("H94").Value = 681.5
$ws.Range("I94").Value = 560.7273
$ws.Range("K94").Value = 560.7273
$ws.Range("M94").Value = -109.7273
$ws.Range("H134").Value = 2461.0334
$ws.Range("I134").Value = 2457.0698
$ws.Range("J134").Value = 2471.0588
$ws.Range("K134").Value = 7371.209400000001
$ws.Range("L134").Value = 7413.176399999999
$ws.Range("M134").Value = -4836.209400000001
$ws.Range("N134").Value = -12483.1764

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13161633
$ws.Range("I58").Value = 3017.6843
$ws.Range("J58").Value = 26320248
$ws.Range("K58").Value = 3017.6843
$ws.Range("L58").Value = 26320248
$ws.Range("M58").Value = -2814.6843
$ws.Range("N58").Value = -26320654
$ws.Range("H132").Value = 2650.2068
$ws.Range("I132").Value = 2811.2856
$ws.Range("J132").Value = 2499.8667
$ws.Range("K132").Value = 8433.856800000001
$ws.Range("L132").Value = 7499.6001
$ws.Range("M132").Value = -5903.856800000001
$ws.Range("N132").Value = -12559.6001
$ws.Range("H134").Value = 14289537
$ws.Range("I134").Value = 21743612
$ws.Range("J134").Value = 2559.5
$ws.Range("K134").Value = 65230836
$ws.Range("L134").Value = 7678.5
$ws.Range("M134").Value = -65228301
$ws.Range("N134").Value = -12748.5
$ws.Range("H136").Value = 13161633
$ws.Range("I136").Value = 3017.6843
$ws.Range("J136").Value = 26320248
$ws.Range("K136").Value = 9053.052899999999
$ws.Range("L136").Value = 78960744
$ws.Range("M136").Value = -6503.052899999999
$ws.Range("N136").Value = -78965844

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 757.125
$ws.Range("I4").Value = 72.25
$ws.Range("J4").Value = 985.4167
$ws.Range("K4").Value = 216.75
$ws.Range("L4").Value = 2956.2501
$ws.Range("M4").Value = -104.75
$ws.Range("N4").Value = -3180.2501
$ws.Range("H57").Value = 3076.25
$ws.Range("I57").Value = 305
$ws.Range("K57").Value = 915
$ws.Range("M57").Value = -356
$ws.Range("H63").Value = 15236.125
$ws.Range("J63").Value = 23198
$ws.Range("L63").Value = 69594
$ws.Range("N63").Value = -71092
$ws.Range("H66").Value = 15236.125
$ws.Range("J66").Value = 23198
$ws.Range("L66").Value = 208782
$ws.Range("N66").Value = -216270
$ws.Range("H113").Value = 751.9583
$ws.Range("I113").Value = 667.6667
$ws.Range("J113").Value = 764
$ws.Range("K113").Value = 2003.0001
$ws.Range("L113").Value = 2292
$ws.Range("M113").Value = 166.9999
$ws.Range("N113").Value = -6632

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 29800
$ws.Range("J75").Value = 29800
$ws.Range("L75").Value = 29800
$ws.Range("N75").Value = -31548
$ws.Range("H78").Value = 29800
$ws.Range("J78").Value = 29800
$ws.Range("L78").Value = 89400
$ws.Range("N78").Value = -98136
$ws.Range("H132").Value = 3528.087
$ws.Range("I132").Value = 3519
$ws.Range("J132").Value = 3539.9
$ws.Range("K132").Value = 10557
$ws.Range("L132").Value = 10619.7
$ws.Range("M132").Value = -8027
$ws.Range("N132").Value = -15679.7

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 20124.5
$ws.Range("J62").Value = 20124.5
$ws.Range("L62").Value = 20124.5
$ws.Range("N62").Value = -21372.5
$ws.Range("H63").Value = 16666.666
$ws.Range("J63").Value = 16666.666
$ws.Range("L63").Value = 16666.666
$ws.Range("N63").Value = -18164.666
$ws.Range("H65").Value = 20124.5
$ws.Range("J65").Value = 20124.5
$ws.Range("L65").Value = 60373.5
$ws.Range("N65").Value = -66613.5
$ws.Range("H66").Value = 16666.666
$ws.Range("J66").Value = 16666.666
$ws.Range("L66").Value = 49999.99800000001
$ws.Range("N66").Value = -57487.99800000001
$ws.Range("H132").Value = 4326.9614
$ws.Range("I132").Value = 1632.2609
$ws.Range("J132").Value = 6464.1377
$ws.Range("K132").Value = 4896.7827
$ws.Range("L132").Value = 19392.4131
$ws.Range("M132").Value = -2366.7827
$ws.Range("N132").Value = -24452.4131
$ws.Range("H136").Value = 2549.0488
$ws.Range("I136").Value = 2491.4644
$ws.Range("J136").Value = 2673.077
$ws.Range("K136").Value = 7474.3932
$ws.Range("L136").Value = 8019.231000000001
$ws.Range("M136").Value = -4924.3932
$ws.Range("N136").Value = -13119.231

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H81").Value = 1709.1
$ws.Range("I81").Value = 1543.5
$ws.Range("J81").Value = 1957.5
$ws.Range("K81").Value = 3087
$ws.Range("L81").Value = 3915
$ws.Range("M81").Value = -2026
$ws.Range("N81").Value = -6037
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H84").Value = 1709.1
$ws.Range("I84").Value = 1543.5
$ws.Range("J84").Value = 1957.5
$ws.Range("K84").Value = 15435
$ws.Range("L84").Value = 19575
$ws.Range("M84").Value = -10131
$ws.Range("N84").Value = -30183
$ws.Range("H107").Value = 283.92856
$ws.Range("I107").Value = 298.1111
$ws.Range("J107").Value = 258.4
$ws.Range("K107").Value = 894.3333
$ws.Range("L107").Value = 775.1999999999999
$ws.Range("M107").Value = 1025.6667
$ws.Range("N107").Value = -4615.2
$ws.Range("H113").Value = 1031.1765
$ws.Range("I113").Value = 533.8333
$ws.Range("J113").Value = 2224.8
$ws.Range("K113").Value = 1601.4999
$ws.Range("L113").Value = 6674.400000000001
$ws.Range("M113").Value = 568.5001
$ws.Range("N113").Value = -11014.4
$ws.Range("H122").Value = 2361.25
$ws.Range("I122").Value = 1521.6154
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 4564.8462
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -2114.8462
$ws.Range("N122").Value = -22898.9995

Write-Host "Applied 241 cell updates across 8 sheets"